# Generate Report for Handback
#
# For the "ce40bf71-6f36-4ab0-a87a-6273e535b378" source file, a new handback
# (target) file has been processed.  This fills in the "Latest Target File",
# "Latest Handback File", "Latest Handback DateTime" and "Error Detail"
# columns (I, J, K, P) for row 7 on both the zh-cn and de-de report sheets,
# and widens the "Error Detail" column so the message is readable.

$wb = $excel.ActiveWorkbook

$hyperlinkUnderlineColor = 15570276   # BGR for RGB(100,149,237) / ARGB FF6495ED, matches the workbook's HyperLink style

function Update-HandbackRow {
    param(
        [string]$SheetName,
        [string]$XlfName,
        [string]$HandbackDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Latest Target File: now available, shown as a hyperlink to the latest
    # revision of the handback markdown file (same target as column A's link).
    $targetCell = $ws.Range("I7")
    $targetCell.Value = "ce40bf71-6f36-4ab0-a87a-6273e535b378.md"
    $ws.Hyperlinks.Add($targetCell, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1ed681b73bd74cf782668cb33d65d32908260d0/e2e/ce40bf71-6f36-4ab0-a87a-6273e535b378.md", "", "", "ce40bf71-6f36-4ab0-a87a-6273e535b378.md") | Out-Null
    $targetCell.Font.Underline = $true
    $targetCell.Font.Color = $hyperlinkUnderlineColor

    # Latest Handback File
    $ws.Range("J7").Value = $XlfName

    # Latest Handback DateTime
    $ws.Range("K7").Value = $HandbackDateTime

    # Error Detail: the handed-back file was generated from an older commit
    # than the current source, so flag the mismatch.
    $ws.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ab9355b0be241d2feb752967aa0d6482c1c81b64/e2e/ce40bf71-6f36-4ab0-a87a-6273e535b378.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1ed681b73bd74cf782668cb33d65d32908260d0/e2e/ce40bf71-6f36-4ab0-a87a-6273e535b378.md."

    # Widen the Error Detail column so the new message is legible.
    $ws.Columns.Item(16).ColumnWidth = 39.15
}

Update-HandbackRow "zh-cn" "ce40bf71-6f36-4ab0-a87a-6273e535b378.3b33b5b9f6585b808bd3d5d8ce82459736b46182.zh-cn.xlf" "2016-08-20 16:52:14"
Update-HandbackRow "de-de" "ce40bf71-6f36-4ab0-a87a-6273e535b378.3b33b5b9f6585b808bd3d5d8ce82459736b46182.de-de.xlf" "2016-08-20 16:52:20"

Write-Host "Handback report updated for ce40bf71-6f36-4ab0-a87a-6273e535b378 on zh-cn and de-de sheets."
